$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15005
$ws1.Range("F3").Value = 18995
$ws1.Range("F14").Value = 153
$ws1.Range("F16").Value = 63
$ws1.Range("F17").Value = 1460
$ws1.Range("F22").Value = 7918
$ws1.Range("F27").Value = 1244
$ws1.Range("F29").Value = 6046
$ws1.Range("F31").Value = 72
$ws1.Range("F32").Value = 169
$ws1.Range("F34").Value = 280
$ws1.Range("F35").Value = 5415
$ws1.Range("F36").Value = 179
$ws1.Range("F37").Value = 9
$ws1.Range("F38").Value = 27

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15005
$ws4.Range("F3").Value = 18995
$ws4.Range("F14").Value = 153
$ws4.Range("F16").Value = 63
$ws4.Range("F17").Value = 1460
$ws4.Range("F23").Value = 7918
$ws4.Range("F28").Value = 1244
$ws4.Range("F32").Value = 6046
$ws4.Range("F34").Value = 72
$ws4.Range("F35").Value = 169
$ws4.Range("F37").Value = 280
$ws4.Range("F38").Value = 5415
$ws4.Range("F39").Value = 180
$ws4.Range("F40").Value = 9
$ws4.Range("F41").Value = 27

Write-Host "Done updating 想去人数 (attendee counts)."
